$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Y and AA (date-text) columns are treated as text to avoid Excel auto-converting
# "YYYY-MM-DD" strings into date serial numbers when the value is assigned.
$ws.Range("Y2:Y22").NumberFormat = "@"
$ws.Range("AA2:AA22").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,1).Value() = 111809318
$ws.Cells.Item(2,2).Value() = 56398
$ws.Cells.Item(2,4).Value() = "NT"
$ws.Cells.Item(2,5).Value() = 100109
$ws.Cells.Item(2,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(2,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(2,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(2,17).Value() = 610793.1422167002
$ws.Cells.Item(2,18).Value() = 7180891.290986502
$ws.Cells.Item(2,25).Value() = "2023-08-31"
$ws.Cells.Item(2,26).Value() = "14:15"
$ws.Cells.Item(2,27).Value() = "2023-08-31"
$ws.Cells.Item(2,28).Value() = "14:15"

# Row 3
$ws.Cells.Item(3,1).Value() = 111808515
$ws.Cells.Item(3,2).Value() = 56398
$ws.Cells.Item(3,4).Value() = "NT"
$ws.Cells.Item(3,5).Value() = 100109
$ws.Cells.Item(3,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(3,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(3,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(3,17).Value() = 610986.3631281323
$ws.Cells.Item(3,18).Value() = 7181120.765008576
$ws.Cells.Item(3,25).Value() = "2023-08-31"
$ws.Cells.Item(3,26).Value() = "13:29"
$ws.Cells.Item(3,27).Value() = "2023-08-31"
$ws.Cells.Item(3,28).Value() = "13:29"

# Row 4
$ws.Cells.Item(4,1).Value() = 111807821
$ws.Cells.Item(4,2).Value() = 56398
$ws.Cells.Item(4,4).Value() = "NT"
$ws.Cells.Item(4,5).Value() = 100109
$ws.Cells.Item(4,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(4,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(4,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(4,17).Value() = 610967.5972068857
$ws.Cells.Item(4,18).Value() = 7181002.477957427
$ws.Cells.Item(4,25).Value() = "2023-08-31"
$ws.Cells.Item(4,26).Value() = "12:51"
$ws.Cells.Item(4,27).Value() = "2023-08-31"
$ws.Cells.Item(4,28).Value() = "12:51"

# Row 5
$ws.Cells.Item(5,1).Value() = 111808676
$ws.Cells.Item(5,2).Value() = 56398
$ws.Cells.Item(5,4).Value() = "NT"
$ws.Cells.Item(5,5).Value() = 100109
$ws.Cells.Item(5,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(5,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(5,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(5,17).Value() = 610906.3771198173
$ws.Cells.Item(5,18).Value() = 7180955.75143602
$ws.Cells.Item(5,25).Value() = "2023-08-31"
$ws.Cells.Item(5,26).Value() = "13:40"
$ws.Cells.Item(5,27).Value() = "2023-08-31"
$ws.Cells.Item(5,28).Value() = "13:40"

# Row 6
$ws.Cells.Item(6,1).Value() = 111808328
$ws.Cells.Item(6,2).Value() = 90087
$ws.Cells.Item(6,4).Value() = "LC"
$ws.Cells.Item(6,5).Value() = 3298
$ws.Cells.Item(6,6).Value() = "Trådticka"
$ws.Cells.Item(6,7).Value() = "Climacocystis borealis"
$ws.Cells.Item(6,8).Value() = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(6,17).Value() = 611008.4619706698
$ws.Cells.Item(6,18).Value() = 7181032.205813259
$ws.Cells.Item(6,25).Value() = "2023-08-31"
$ws.Cells.Item(6,26).Value() = "13:23"
$ws.Cells.Item(6,27).Value() = "2023-08-31"
$ws.Cells.Item(6,28).Value() = "13:23"

# Row 7
$ws.Cells.Item(7,1).Value() = 111808387
$ws.Cells.Item(7,2).Value() = 89419
$ws.Cells.Item(7,4).Value() = "NT"
$ws.Cells.Item(7,5).Value() = 1204
$ws.Cells.Item(7,6).Value() = "Gränsticka"
$ws.Cells.Item(7,7).Value() = "Phellopilus nigrolimitatus"
$ws.Cells.Item(7,8).Value() = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Cells.Item(7,17).Value() = 611016.8359391808
$ws.Cells.Item(7,18).Value() = 7181141.984797659
$ws.Cells.Item(7,25).Value() = "2023-08-31"
$ws.Cells.Item(7,26).Value() = "13:26"
$ws.Cells.Item(7,27).Value() = "2023-08-31"
$ws.Cells.Item(7,28).Value() = "13:26"

# Row 8
$ws.Cells.Item(8,1).Value() = 111807055
$ws.Cells.Item(8,2).Value() = 77515
$ws.Cells.Item(8,4).Value() = "NT"
$ws.Cells.Item(8,5).Value() = 6425
$ws.Cells.Item(8,6).Value() = "Garnlav"
$ws.Cells.Item(8,7).Value() = "Alectoria sarmentosa"
$ws.Cells.Item(8,8).Value() = "(Ach.) Ach."
$ws.Cells.Item(8,17).Value() = 610666.4119294117
$ws.Cells.Item(8,18).Value() = 7181042.722880279
$ws.Cells.Item(8,25).Value() = "2023-08-31"
$ws.Cells.Item(8,26).Value() = "12:10"
$ws.Cells.Item(8,27).Value() = "2023-08-31"
$ws.Cells.Item(8,28).Value() = "12:10"

# Row 9
$ws.Cells.Item(9,1).Value() = 111808022
$ws.Cells.Item(9,2).Value() = 90087
$ws.Cells.Item(9,4).Value() = "LC"
$ws.Cells.Item(9,5).Value() = 3298
$ws.Cells.Item(9,6).Value() = "Trådticka"
$ws.Cells.Item(9,7).Value() = "Climacocystis borealis"
$ws.Cells.Item(9,8).Value() = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(9,17).Value() = 611005.4375152331
$ws.Cells.Item(9,18).Value() = 7181032.949711116
$ws.Cells.Item(9,25).Value() = "2023-08-31"
$ws.Cells.Item(9,26).Value() = "12:58"
$ws.Cells.Item(9,27).Value() = "2023-08-31"
$ws.Cells.Item(9,28).Value() = "12:58"

# Row 10
$ws.Cells.Item(10,1).Value() = 111807814
$ws.Cells.Item(10,2).Value() = 56398
$ws.Cells.Item(10,4).Value() = "NT"
$ws.Cells.Item(10,5).Value() = 100109
$ws.Cells.Item(10,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(10,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(10,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(10,17).Value() = 610955.7664983921
$ws.Cells.Item(10,18).Value() = 7180998.193866283
$ws.Cells.Item(10,25).Value() = "2023-08-31"
$ws.Cells.Item(10,26).Value() = "12:50"
$ws.Cells.Item(10,27).Value() = "2023-08-31"
$ws.Cells.Item(10,28).Value() = "12:50"

# Row 11
$ws.Cells.Item(11,1).Value() = 111809897
$ws.Cells.Item(11,2).Value() = 85715
$ws.Cells.Item(11,4).Value() = "NT"
$ws.Cells.Item(11,5).Value() = 510
$ws.Cells.Item(11,6).Value() = "Doftskinn"
$ws.Cells.Item(11,7).Value() = "Cystostereum murrayi"
$ws.Cells.Item(11,8).Value() = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Cells.Item(11,17).Value() = 610718.603132805
$ws.Cells.Item(11,18).Value() = 7180857.334717941
$ws.Cells.Item(11,25).Value() = "2023-08-31"
$ws.Cells.Item(11,26).Value() = "14:50"
$ws.Cells.Item(11,27).Value() = "2023-08-31"
$ws.Cells.Item(11,28).Value() = "14:50"

# Row 12
$ws.Cells.Item(12,1).Value() = 111807370
$ws.Cells.Item(12,2).Value() = 56398
$ws.Cells.Item(12,4).Value() = "NT"
$ws.Cells.Item(12,5).Value() = 100109
$ws.Cells.Item(12,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(12,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(12,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(12,17).Value() = 610667.3767981895
$ws.Cells.Item(12,18).Value() = 7181039.764941735
$ws.Cells.Item(12,25).Value() = "2023-08-31"
$ws.Cells.Item(12,26).Value() = "12:36"
$ws.Cells.Item(12,27).Value() = "2023-08-31"
$ws.Cells.Item(12,28).Value() = "12:36"

# Row 13
$ws.Cells.Item(13,1).Value() = 111808000
$ws.Cells.Item(13,2).Value() = 56398
$ws.Cells.Item(13,4).Value() = "NT"
$ws.Cells.Item(13,5).Value() = 100109
$ws.Cells.Item(13,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(13,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(13,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(13,17).Value() = 610972.8318012832
$ws.Cells.Item(13,18).Value() = 7180976.585010829
$ws.Cells.Item(13,25).Value() = "2023-08-31"
$ws.Cells.Item(13,26).Value() = "12:56"
$ws.Cells.Item(13,27).Value() = "2023-08-31"
$ws.Cells.Item(13,28).Value() = "12:56"

# Row 14
$ws.Cells.Item(14,1).Value() = 111806969
$ws.Cells.Item(14,2).Value() = 56398
$ws.Cells.Item(14,4).Value() = "NT"
$ws.Cells.Item(14,5).Value() = 100109
$ws.Cells.Item(14,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(14,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(14,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(14,17).Value() = 610695.5210812307
$ws.Cells.Item(14,18).Value() = 7181007.871029559
$ws.Cells.Item(14,25).Value() = "2023-08-31"
$ws.Cells.Item(14,26).Value() = "12:09"
$ws.Cells.Item(14,27).Value() = "2023-08-31"
$ws.Cells.Item(14,28).Value() = "12:09"

# Row 15
$ws.Cells.Item(15,1).Value() = 111809026
$ws.Cells.Item(15,2).Value() = 56398
$ws.Cells.Item(15,4).Value() = "NT"
$ws.Cells.Item(15,5).Value() = 100109
$ws.Cells.Item(15,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(15,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(15,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(15,17).Value() = 610772.3941544103
$ws.Cells.Item(15,18).Value() = 7180884.969268824
$ws.Cells.Item(15,25).Value() = "2023-08-31"
$ws.Cells.Item(15,26).Value() = "13:56"
$ws.Cells.Item(15,27).Value() = "2023-08-31"
$ws.Cells.Item(15,28).Value() = "13:56"

# Row 16
$ws.Cells.Item(16,1).Value() = 111808957
$ws.Cells.Item(16,2).Value() = 56398
$ws.Cells.Item(16,4).Value() = "NT"
$ws.Cells.Item(16,5).Value() = 100109
$ws.Cells.Item(16,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(16,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(16,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(16,17).Value() = 610787.509024983
$ws.Cells.Item(16,18).Value() = 7180858.155172224
$ws.Cells.Item(16,25).Value() = "2023-08-31"
$ws.Cells.Item(16,26).Value() = "13:52"
$ws.Cells.Item(16,27).Value() = "2023-08-31"
$ws.Cells.Item(16,28).Value() = "13:52"

# Row 17
$ws.Cells.Item(17,1).Value() = 111808659
$ws.Cells.Item(17,2).Value() = 56398
$ws.Cells.Item(17,4).Value() = "NT"
$ws.Cells.Item(17,5).Value() = 100109
$ws.Cells.Item(17,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(17,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(17,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(17,17).Value() = 610912.6426496292
$ws.Cells.Item(17,18).Value() = 7180936.738188162
$ws.Cells.Item(17,25).Value() = "2023-08-31"
$ws.Cells.Item(17,26).Value() = "13:39"
$ws.Cells.Item(17,27).Value() = "2023-08-31"
$ws.Cells.Item(17,28).Value() = "13:39"

# Row 18
$ws.Cells.Item(18,1).Value() = 111807777
$ws.Cells.Item(18,2).Value() = 56398
$ws.Cells.Item(18,4).Value() = "NT"
$ws.Cells.Item(18,5).Value() = 100109
$ws.Cells.Item(18,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(18,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(18,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(18,17).Value() = 610955.5779051523
$ws.Cells.Item(18,18).Value() = 7181003.318500374
$ws.Cells.Item(18,25).Value() = "2023-08-31"
$ws.Cells.Item(18,26).Value() = "12:48"
$ws.Cells.Item(18,27).Value() = "2023-08-31"
$ws.Cells.Item(18,28).Value() = "12:48"

# Row 19
$ws.Cells.Item(19,1).Value() = 111806888
$ws.Cells.Item(19,2).Value() = 56398
$ws.Cells.Item(19,4).Value() = "NT"
$ws.Cells.Item(19,5).Value() = 100109
$ws.Cells.Item(19,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(19,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(19,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(19,17).Value() = 610667.2784075569
$ws.Cells.Item(19,18).Value() = 7181030.781246905
$ws.Cells.Item(19,25).Value() = "1970-01-01"
$ws.Cells.Item(19,26).Value() = "00:00"
$ws.Cells.Item(19,27).Value() = "1970-01-01"
$ws.Cells.Item(19,28).Value() = "00:00"

# Row 20
$ws.Cells.Item(20,1).Value() = 111808532
$ws.Cells.Item(20,2).Value() = 89423
$ws.Cells.Item(20,4).Value() = "NT"
$ws.Cells.Item(20,5).Value() = 5432
$ws.Cells.Item(20,6).Value() = "Granticka"
$ws.Cells.Item(20,7).Value() = "Porodaedalea chrysoloma"
$ws.Cells.Item(20,8).Value() = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(20,17).Value() = 610963.6803355663
$ws.Cells.Item(20,18).Value() = 7181097.266412931
$ws.Cells.Item(20,25).Value() = "2023-08-31"
$ws.Cells.Item(20,26).Value() = "13:32"
$ws.Cells.Item(20,27).Value() = "2023-08-31"
$ws.Cells.Item(20,28).Value() = "13:32"

# Row 21
$ws.Cells.Item(21,1).Value() = 111807092
$ws.Cells.Item(21,2).Value() = 77515
$ws.Cells.Item(21,4).Value() = "NT"
$ws.Cells.Item(21,5).Value() = 6425
$ws.Cells.Item(21,6).Value() = "Garnlav"
$ws.Cells.Item(21,7).Value() = "Alectoria sarmentosa"
$ws.Cells.Item(21,8).Value() = "(Ach.) Ach."
$ws.Cells.Item(21,17).Value() = 610678.3995887624
$ws.Cells.Item(21,18).Value() = 7181042.735301252
$ws.Cells.Item(21,25).Value() = "2023-08-31"
$ws.Cells.Item(21,26).Value() = "12:11"
$ws.Cells.Item(21,27).Value() = "2023-08-31"
$ws.Cells.Item(21,28).Value() = "12:11"

# Row 22
$ws.Cells.Item(22,1).Value() = 111808817
$ws.Cells.Item(22,2).Value() = 56398
$ws.Cells.Item(22,4).Value() = "NT"
$ws.Cells.Item(22,5).Value() = 100109
$ws.Cells.Item(22,6).Value() = "Tretåig hackspett"
$ws.Cells.Item(22,7).Value() = "Picoides tridactylus"
$ws.Cells.Item(22,8).Value() = "(Linnaeus, 1758)"
$ws.Cells.Item(22,17).Value() = 610921.7319367616
$ws.Cells.Item(22,18).Value() = 7180934.079081071
$ws.Cells.Item(22,25).Value() = "2023-08-31"
$ws.Cells.Item(22,26).Value() = "13:45"
$ws.Cells.Item(22,27).Value() = "2023-08-31"
$ws.Cells.Item(22,28).Value() = "13:45"

# Restore General number format for the date-text columns while keeping their text values intact
$ws.Range("Y2:Y22").ClearFormats()
$ws.Range("AA2:AA22").ClearFormats()

# Row 11 now carries what used to be row 16's content, which has extra (empty) J, N and AF cells.
$ws.Cells.Item(11,10).Value() = ""
$ws.Cells.Item(11,14).Value() = ""
$ws.Cells.Item(11,32).Value() = ""

# Row 16 now carries what used to be row 18's content, which does NOT have J, N, AF cells.
$ws.Cells.Item(16,10).ClearContents()
$ws.Cells.Item(16,14).ClearContents()
$ws.Cells.Item(16,32).ClearContents()